$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 7752
$ws.Range("F8").Value = 2111
$ws.Range("F9").Value = 8576
$ws.Range("F12").Value = 86
$ws.Range("F13").Value = 5723
$ws.Range("F15").Value = 2674
$ws.Range("F16").Value = 1179
$ws.Range("F22").Value = 571
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 3732
$ws.Range("F25").Value = 74
$ws.Range("F29").Value = 3233
$ws.Range("F30").Value = 59
$ws.Range("F34").Value = 350
$ws.Range("F35").Value = 1133
$ws.Range("F36").Value = 685
$ws.Range("F39").Value = 2829
$ws.Range("G39").Value = 5
$ws.Range("F40").Value = 57
$ws.Range("F43").Value = 3258
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7752
$ws.Range("F8").Value = 2111
$ws.Range("F9").Value = 8576
$ws.Range("F11").Value = 86
$ws.Range("F12").Value = 5723
$ws.Range("F14").Value = 2674
$ws.Range("F15").Value = 1179
$ws.Range("F22").Value = 571
$ws.Range("F24").Value = 3732
$ws.Range("F25").Value = 74
$ws.Range("F29").Value = 3233
$ws.Range("F30").Value = 59
$ws.Range("F33").Value = 350
$ws.Range("F35").Value = 1133
$ws.Range("F36").Value = 685
$ws.Range("F40").Value = 2829
$ws.Range("G40").Value = 5
$ws.Range("F41").Value = 57
$ws.Range("F44").Value = 3258
